$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 421.25
$ws.Range("I11").Value = 421.25
$ws.Range("K11").Value = 421.25
$ws.Range("M11").Value = -281.25
$ws.Range("H43").Value = 1997.5
$ws.Range("I43").Value = 1997.5
$ws.Range("K43").Value = 1997.5
$ws.Range("M43").Value = -1928.5
$ws.Range("H69").Value = 13000
$ws.Range("I69").Value = 13000
$ws.Range("K69").Value = 39000
$ws.Range("M69").Value = -38126
$ws.Range("H72").Value = 13000
$ws.Range("I72").Value = 13000
$ws.Range("K72").Value = 117000
$ws.Range("M72").Value = -112632
$ws.Range("H76").Value = 5494.4165
$ws.Range("I76").Value = 2500
$ws.Range("J76").Value = 6093.3
$ws.Range("K76").Value = 2500
$ws.Range("L76").Value = 6093.3
$ws.Range("M76").Value = -2185
$ws.Range("N76").Value = -6723.3
$ws.Range("H79").Value = 5494.4165
$ws.Range("I79").Value = 2500
$ws.Range("J79").Value = 6093.3
$ws.Range("K79").Value = 2500
$ws.Range("L79").Value = 6093.3
$ws.Range("M79").Value = -1408
$ws.Range("N79").Value = -8277.299999999999
$ws.Range("H88").Value = 1724.0769
$ws.Range("J88").Value = 1944.1
$ws.Range("L88").Value = 1944.1
$ws.Range("N88").Value = -2756.1
$ws.Range("H91").Value = 1724.0769
$ws.Range("J91").Value = 1944.1
$ws.Range("L91").Value = 1944.1
$ws.Range("N91").Value = -4752.1
$ws.Range("H113").Value = 5960.5
$ws.Range("I113").Value = 5153
$ws.Range("K113").Value = 5153
$ws.Range("M113").Value = -1899
$ws.Range("H132").Value = 3330
$ws.Range("I132").Value = 3274.7058
$ws.Range("K132").Value = 9824.117400000001
$ws.Range("M132").Value = -7294.117400000001
$ws.Range("H137").Value = 8218.182000000001
$ws.Range("I137").Value = 6540
$ws.Range("K137").Value = 19620
$ws.Range("M137").Value = -17070
$ws.Range("H141").Value = 800
$ws.Range("I141").Value = 825
$ws.Range("J141").Value = 750
$ws.Range("K141").Value = 2475
$ws.Range("L141").Value = 2250
$ws.Range("M141").Value = 2705
$ws.Range("N141").Value = -12610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2801
$ws.Range("I2").Value = 2126.25
$ws.Range("K2").Value = 2126.25
$ws.Range("M2").Value = -2013.25
$ws.Range("H32").Value = 20453.334
$ws.Range("I32").Value = 12597.685
$ws.Range("K32").Value = 12597.685
$ws.Range("M32").Value = -12310.685
$ws.Range("H45").Value = 2570.9375
$ws.Range("J45").Value = 3571.4285
$ws.Range("L45").Value = 3571.4285
$ws.Range("N45").Value = -4325.4285
$ws.Range("H110").Value = 10999.667
$ws.Range("I110").Value = 10999.667
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 10999.667
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -8954.666999999999
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 2801
$ws.Range("I116").Value = 2126.25
$ws.Range("K116").Value = 2126.25
$ws.Range("M116").Value = 167.75
$ws.Range("H122").Value = 1254002.1
$ws.Range("I122").Value = 2003601.6
$ws.Range("J122").Value = 4669.6665
$ws.Range("K122").Value = 6010804.800000001
$ws.Range("L122").Value = 14008.9995
$ws.Range("M122").Value = -6008354.800000001
$ws.Range("N122").Value = -18908.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2801
$ws.Range("I3").Value = 2126.25
$ws.Range("K3").Value = 2126.25
$ws.Range("M3").Value = -2012.25
$ws.Range("H94").Value = 2507.6667
$ws.Range("I94").Value = 2508
$ws.Range("J94").Value = 2507
$ws.Range("K94").Value = 2508
$ws.Range("L94").Value = 2507
$ws.Range("M94").Value = -2057
$ws.Range("N94").Value = -3409
$ws.Range("H134").Value = 3169.5908
$ws.Range("I134").Value = 1579.0834
$ws.Range("K134").Value = 4737.2502
$ws.Range("M134").Value = -2202.2502

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 125000
$ws.Range("J88").Value = 125000
$ws.Range("L88").Value = 125000
$ws.Range("N88").Value = -125812
$ws.Range("H91").Value = 125000
$ws.Range("J91").Value = 125000
$ws.Range("L91").Value = 125000
$ws.Range("N91").Value = -127808
$ws.Range("H107").Value = 398
$ws.Range("I107").Value = 352
$ws.Range("J107").Value = 417.7143
$ws.Range("K107").Value = 352
$ws.Range("L107").Value = 417.7143
$ws.Range("M107").Value = 1568
$ws.Range("N107").Value = -4257.7143
$ws.Range("H132").Value = 2261.879
$ws.Range("I132").Value = 2165.077
$ws.Range("K132").Value = 6495.231000000001
$ws.Range("M132").Value = -3965.231000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9434.9
$ws.Range("J80").Value = 9431.875
$ws.Range("L80").Value = 9431.875
$ws.Range("N80").Value = -11427.875
$ws.Range("H83").Value = 9434.9
$ws.Range("J83").Value = 9431.875
$ws.Range("L83").Value = 47159.375
$ws.Range("N83").Value = -57143.375
$ws.Range("H102").Value = 4184.5
$ws.Range("I102").Value = 2825.3333
$ws.Range("K102").Value = 2825.3333
$ws.Range("M102").Value = -1203.3333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4840.8335
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H61").Value = 5999.4585
$ws.Range("I61").Value = 5999.4287
$ws.Range("K61").Value = 5999.4287
$ws.Range("M61").Value = -5797.4287
$ws.Range("H82").Value = 2168.9333
$ws.Range("J82").Value = 1520.7142
$ws.Range("L82").Value = 1520.7142
$ws.Range("N82").Value = -2242.7142
$ws.Range("H85").Value = 2168.9333
$ws.Range("J85").Value = 1520.7142
$ws.Range("L85").Value = 1520.7142
$ws.Range("N85").Value = -4016.7142
$ws.Range("H100").Value = 4643.2856
$ws.Range("I100").Value = 4583.8335
$ws.Range("K100").Value = 4583.8335
$ws.Range("M100").Value = -4042.8335
$ws.Range("H113").Value = 5999.4585
$ws.Range("I113").Value = 5999.4287
$ws.Range("K113").Value = 5999.4287
$ws.Range("M113").Value = -3829.4287
$ws.Range("H122").Value = 8699.875
$ws.Range("I122").Value = 7399.75
$ws.Range("K122").Value = 22199.25
$ws.Range("M122").Value = -19749.25
$ws.Range("H126").Value = 4840.8335
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4739.7856
$ws.Range("I132").Value = 2392.8333
$ws.Range("K132").Value = 7178.499899999999
$ws.Range("M132").Value = -4648.499899999999
